$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is being trimmed down to a single-year export: drop the blank
# spacer row that used to sit between the title block and the "Area"
# sub-header (it only ever held one empty, styled cell) ...
$ws.Range("A3").EntireRow.Delete()

# ... and drop the 1989 / 2002 columns, keeping only the 2014 figure
# (it shifts left into column B).
$ws.Range("B1:C1").EntireColumn.Delete()

# The "(according to the population census data)" note is removed from
# the export entirely, leaving row 2 blank.
$ws.Range("A2:B2").Clear()

# The title row's trailing empty styled cell is removed too.
$ws.Range("B1").Clear()

# Every row gets the new, taller row height used by the refreshed layout
# (this also appends a handful of blank rows below the table).
$ws.Range("A1:B11").RowHeight = 20.1
